# Scheduled market-price refresh: recompute currentAveragePrice* / LevePrice* / LeveProfit*
# columns (H:N) per leve row from refreshed Universalis data, across all eight crafter sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17: One for the Road
$ws.Range("H17").Value = 1662.6216
$ws.Range("J17").Value = 1662.6216
$ws.Range("L17").Value = 4987.864799999999
$ws.Range("N17").Value = -5323.864799999999
# Row 28: The Writing Is Not on the Wall
$ws.Range("H28").Value = 80361.69500000001
$ws.Range("I28").Value = 758.1429000000001
$ws.Range("J28").Value = 173232.5
$ws.Range("K28").Value = 758.1429000000001
$ws.Range("L28").Value = 173232.5
$ws.Range("M28").Value = -273.1429000000001
$ws.Range("N28").Value = -174202.5
# Row 55: A Real Smooth Move
$ws.Range("H55").Value = 145
$ws.Range("I55").Value = 110.833336
$ws.Range("K55").Value = 110.833336
$ws.Range("M55").Value = 103.166664
# Row 70: Consecrating Congregation
$ws.Range("H70").Value = 1469.9445
$ws.Range("I70").Value = 1230.6
$ws.Range("J70").Value = 2666.6667
$ws.Range("K70").Value = 3691.8
$ws.Range("L70").Value = 8000.000100000001
$ws.Range("M70").Value = -3421.8
$ws.Range("N70").Value = -8540.000100000001
# Row 73: Curbing the Contagion (L)
$ws.Range("H73").Value = 1469.9445
$ws.Range("I73").Value = 1230.6
$ws.Range("J73").Value = 2666.6667
$ws.Range("K73").Value = 3691.8
$ws.Range("L73").Value = 8000.000100000001
$ws.Range("M73").Value = -2755.8
$ws.Range("N73").Value = -9872.000100000001
# Row 93: Spellbound
$ws.Range("H93").Value = 55066.332
$ws.Range("J93").Value = 55066.332
$ws.Range("L93").Value = 55066.332
$ws.Range("N93").Value = -60058.332
# Row 103: Let Loose the Juice
$ws.Range("H103").Value = 205.23529
$ws.Range("I103").Value = 133.5
$ws.Range("K103").Value = 400.5
$ws.Range("M103").Value = 185.5
# Row 106: Making Your Mark
$ws.Range("H106").Value = 62861100
$ws.Range("I106").Value = 73337110
$ws.Range("K106").Value = 73337110
$ws.Range("M106").Value = -73336479
# Row 116: Growing Up
$ws.Range("H116").Value = 8665.6
$ws.Range("I116").Value = 7648.5
$ws.Range("J116").Value = 9828
$ws.Range("K116").Value = 7648.5
$ws.Range("L116").Value = 9828
$ws.Range("M116").Value = -4206.5
$ws.Range("N116").Value = -16712
# Row 127: Liquid Competence
$ws.Range("H127").Value = 3585.5
$ws.Range("I127").Value = 3418.6667
$ws.Range("K127").Value = 10256.0001
$ws.Range("M127").Value = -5296.000100000001
# Row 129: Practical Command
$ws.Range("H129").Value = 12388
$ws.Range("I129").Value = 12866
$ws.Range("K129").Value = 38598
$ws.Range("M129").Value = -33598
# Row 131: Mindful Study
$ws.Range("H131").Value = 797.4
$ws.Range("I131").Value = 797.4
$ws.Range("K131").Value = 2392.2
$ws.Range("M131").Value = 2647.8
# Row 141: Remedy for Reason
$ws.Range("H141").Value = 5927.8184
$ws.Range("I141").Value = 2534.5
$ws.Range("J141").Value = 9999.799999999999
$ws.Range("K141").Value = 7603.5
$ws.Range("L141").Value = 29999.4
$ws.Range("M141").Value = -2423.5
$ws.Range("N141").Value = -40359.39999999999

$ws = $wb.Worksheets.Item("ARM")
# Row 34: Insistent Sallets
$ws.Range("H34").Value = 71250
$ws.Range("I34").Value = 45000
$ws.Range("J34").Value = 80000
$ws.Range("K34").Value = 45000
$ws.Range("L34").Value = 80000
$ws.Range("M34").Value = -44729
$ws.Range("N34").Value = -80542
# Row 64: Don't Scuttle with Scuta
$ws.Range("H64").Value = 37263
$ws.Range("I64").Value = 35000
$ws.Range("J64").Value = 38394.5
$ws.Range("K64").Value = 35000
$ws.Range("L64").Value = 38394.5
$ws.Range("M64").Value = -34752
$ws.Range("N64").Value = -38890.5
# Row 67: Shielded by Bureaucracy (L)
$ws.Range("H67").Value = 37263
$ws.Range("I67").Value = 35000
$ws.Range("J67").Value = 38394.5
$ws.Range("K67").Value = 35000
$ws.Range("L67").Value = 38394.5
$ws.Range("M67").Value = -34142
$ws.Range("N67").Value = -40110.5
# Row 74: As the Bolt Flies
$ws.Range("H74").Value = 5977.222
$ws.Range("I74").Value = 2356.4
$ws.Range("K74").Value = 2356.4
$ws.Range("M74").Value = -1482.4
# Row 77: Heavy Metal Banned (L)
$ws.Range("H77").Value = 5977.222
$ws.Range("I77").Value = 2356.4
$ws.Range("K77").Value = 11782
$ws.Range("M77").Value = -7414
# Row 122: Haste for High Durium
$ws.Range("H122").Value = 2058.862
$ws.Range("I122").Value = 2203.158
$ws.Range("J122").Value = 1784.7
$ws.Range("K122").Value = 6609.474
$ws.Range("L122").Value = 5354.1
$ws.Range("M122").Value = -4159.474
$ws.Range("N122").Value = -10254.1

$ws = $wb.Worksheets.Item("BSM")
# Row 60: Talon Terrors
$ws.Range("H60").Value = 50929
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 50929
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 50929
$ws.Range("M60").ClearContents()
$ws.Range("N60").Value = -52127
# Row 62: Barring the Gates to Foundation
$ws.Range("H62").Value = 40000
$ws.Range("J62").Value = 40000
$ws.Range("L62").Value = 40000
$ws.Range("N62").Value = -41372
# Row 65: Starting Young (L)
$ws.Range("H65").Value = 40000
$ws.Range("J65").Value = 40000
$ws.Range("L65").Value = 120000
$ws.Range("N65").Value = -126864
# Row 86: Through Thick and Thin
$ws.Range("H86").Value = 3995.3333
$ws.Range("I86").Value = 3995.3333
$ws.Range("K86").Value = 3995.3333
$ws.Range("M86").Value = -2872.3333
# Row 89: Piercing Eyes Deserve Piercing Shafts (L)
$ws.Range("H89").Value = 3995.3333
$ws.Range("I89").Value = 3995.3333
$ws.Range("K89").Value = 19976.6665
$ws.Range("M89").Value = -14360.6665
# Row 94: High Steal
$ws.Range("H94").Value = 824.3333
$ws.Range("I94").Value = 845
$ws.Range("J94").Value = 721
$ws.Range("K94").Value = 845
$ws.Range("L94").Value = 721
$ws.Range("M94").Value = -394
$ws.Range("N94").Value = -1623
# Row 107: The Gold Experience
$ws.Range("H107").Value = 3932.15
$ws.Range("I107").Value = 3956.7856
$ws.Range("K107").Value = 3956.7856
$ws.Range("M107").Value = -2036.7856

$ws = $wb.Worksheets.Item("CRP")
# Row 53: A Winning Combo
$ws.Range("H53").Value = 50000
$ws.Range("J53").Value = 50000
$ws.Range("L53").Value = 50000
$ws.Range("N53").Value = -51214
# Row 105: Zelkova, My Love
$ws.Range("H105").Value = 4451.75
$ws.Range("J105").Value = 4998
$ws.Range("L105").Value = 4998
$ws.Range("N105").Value = -8492

$ws = $wb.Worksheets.Item("CUL")
# Row 23: Sweet Smell of Success
$ws.Range("H23").Value = 131.16667
$ws.Range("I23").Value = 66
$ws.Range("J23").Value = 196.33333
$ws.Range("K23").Value = 198
$ws.Range("L23").Value = 588.99999
$ws.Range("M23").Value = 37
$ws.Range("N23").Value = -1058.99999
# Row 39: Bloody Good Tart, This
$ws.Range("H39").Value = 8712.857
$ws.Range("J39").Value = 11340
$ws.Range("L39").Value = 34020
$ws.Range("N39").Value = -34608
# Row 44: No More Dumpster Diving
$ws.Range("H44").Value = 1162.3334
$ws.Range("I44").Value = 239
$ws.Range("J44").Value = 1624
$ws.Range("K44").Value = 717
$ws.Range("L44").Value = 4872
$ws.Range("M44").Value = -319
$ws.Range("N44").Value = -5668
# Row 48: Rise and Dine
$ws.Range("H48").Value = 95
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").ClearContents()
# Row 51: The Perks of Life at Sea
$ws.Range("H51").Value = 190
$ws.Range("I51").Value = 190
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 570
$ws.Range("L51").Value = 0
$ws.Range("M51").Value = -110
$ws.Range("N51").ClearContents()
# Row 107: Slippery Service
$ws.Range("H107").Value = 956.6667
$ws.Range("J107").Value = 956.6667
$ws.Range("L107").Value = 2870.0001
$ws.Range("N107").Value = -6710.0001
# Row 113: Can't Eat Just One
$ws.Range("H113").Value = 2357.4167
$ws.Range("I113").Value = 3932.6667
$ws.Range("J113").Value = 1832.3334
$ws.Range("K113").Value = 11798.0001
$ws.Range("L113").Value = 5497.0002
$ws.Range("M113").Value = -9628.000100000001
$ws.Range("N113").Value = -9837.0002
# Row 129: Comfort Food
$ws.Range("H129").Value = 8784217
$ws.Range("J129").Value = 15162414
$ws.Range("L129").Value = 45487242
$ws.Range("N129").Value = -45497242
# Row 131: The Mountain Steeped
$ws.Range("H131").Value = 2079.5217
$ws.Range("J131").Value = 5289.857
$ws.Range("L131").Value = 15869.571
$ws.Range("N131").Value = -25949.571
# Row 139: Najoothie
$ws.Range("H139").Value = 1972.3334
$ws.Range("I139").Value = 1445.1666
$ws.Range("K139").Value = 4335.4998
$ws.Range("M139").Value = 804.5002000000004
# Row 140: Sweet, Sweet Bean Juice
$ws.Range("H140").Value = 1866.4706
$ws.Range("I140").Value = 1940.875
$ws.Range("J140").Value = 1800.3334
$ws.Range("K140").Value = 5822.625
$ws.Range("L140").Value = 5401.0002
$ws.Range("M140").Value = -642.625
$ws.Range("N140").Value = -15761.0002

$ws = $wb.Worksheets.Item("GSM")
# Row 102: Put the Metal to the Peddle
$ws.Range("H102").Value = 2840.8333
$ws.Range("I102").Value = 2230.5757
$ws.Range("J102").Value = 5078.4443
$ws.Range("K102").Value = 2230.5757
$ws.Range("L102").Value = 5078.4443
$ws.Range("M102").Value = -608.5756999999999
$ws.Range("N102").Value = -8322.444299999999
# Row 132: On Board for Lar
$ws.Range("H132").Value = 3119.5
$ws.Range("I132").Value = 2098.4688
$ws.Range("K132").Value = 6295.4064
$ws.Range("M132").Value = -3765.4064

$ws = $wb.Worksheets.Item("LTW")
# Row 22: Skin off Their Backs
$ws.Range("H22").Value = 805.2778
$ws.Range("I22").Value = 499.86667
$ws.Range("K22").Value = 499.86667
$ws.Range("M22").Value = -204.86667
# Row 27: Fire and Hide
$ws.Range("H27").Value = 805.2778
$ws.Range("I27").Value = 499.86667
$ws.Range("K27").Value = 499.86667
$ws.Range("M27").Value = -392.86667
# Row 100: Tiger in the Sack
$ws.Range("H100").Value = 2046.7142
$ws.Range("I100").Value = 1887.8334
$ws.Range("K100").Value = 1887.8334
$ws.Range("M100").Value = -1346.8334
# Row 122: Hell on Leather
$ws.Range("H122").Value = 4291.577
$ws.Range("I122").Value = 3675.6667
$ws.Range("J122").Value = 6878.4
$ws.Range("K122").Value = 11027.0001
$ws.Range("L122").Value = 20635.2
$ws.Range("M122").Value = -8577.000100000001
$ws.Range("N122").Value = -25535.2

